$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper scratch cell (outside the used A1:E51 range) used to write literal
# text into cells without Excel's automatic "looks like a number" coercion:
# format the scratch cell as Text, set the literal value there, copy it, then
# PasteSpecial (values only) into the destination cell so the destinations
# own (default/General) number format and style are left untouched.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-LiteralText($cellRef, $text) {
    $scratch = $ws.Range("Z1")
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
}

# --- Row 42 / Row 43: WhiteBITCoin moves up to row 42, BabyDogeCoin moves to row 43 ---
Set-LiteralText "B42" "WhiteBITCoin"
Set-LiteralText "C42" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-LiteralText "D42" "17.82"
Set-LiteralText "E42" "  +0.37%  "
Set-LiteralText "B43" "BabyDogeCoin"
Set-LiteralText "C43" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-LiteralText "D43" "0.0₆0333"
Set-LiteralText "E43" "  +3.38%  "

# --- Remaining price / volume (%) updates ---
Set-LiteralText "D2" "68.190.01"
Set-LiteralText "E2" "  -0.82%  "
Set-LiteralText "D3" "2.642.65"
Set-LiteralText "E3" "  -0.59%  "
Set-LiteralText "E4" "  -0.06%  "
Set-LiteralText "D5" "597.80"
Set-LiteralText "E5" "  -0.49%  "
Set-LiteralText "D6" "156.44"
Set-LiteralText "E6" "  +0.55%  "
Set-LiteralText "E7" "  -0.03%  "
Set-LiteralText "E8" "  -0.88%  "
Set-LiteralText "D9" "0.141"
Set-LiteralText "E9" "  +1.63%  "
Set-LiteralText "E10" "  -1.32%  "
Set-LiteralText "E11" "  +0.23%  "
Set-LiteralText "D12" "0.351"
Set-LiteralText "E12" "  +0.34%  "
Set-LiteralText "D13" "27.99"
Set-LiteralText "E13" "  -0.17%  "
Set-LiteralText "E14" "  +0.84%  "
Set-LiteralText "D15" "3.123.02"
Set-LiteralText "E15" "  -0.65%  "
Set-LiteralText "D16" "68.235.55"
Set-LiteralText "E16" "  -0.60%  "
Set-LiteralText "D17" "2.655.72"
Set-LiteralText "E17" "  -0.27%  "
Set-LiteralText "D18" "11.40"
Set-LiteralText "E18" "  -0.46%  "
Set-LiteralText "D19" "363.59"
Set-LiteralText "E19" "  -0.95%  "
Set-LiteralText "E20" "  -1.47%  "
Set-LiteralText "E21" "  +3.28%  "
Set-LiteralText "E22" "  -1.93%  "
Set-LiteralText "E23" "  -3.18%  "
Set-LiteralText "D24" "75.49"
Set-LiteralText "E24" "  +3.85%  "
Set-LiteralText "E25" "  -0.09%  "
Set-LiteralText "D26" "9.76"
Set-LiteralText "E26" "  -2.73%  "
Set-LiteralText "D27" "1.04"
Set-LiteralText "E27" "  +3.36%  "
Set-LiteralText "D28" "2.776.00"
Set-LiteralText "E28" "  -0.98%  "
Set-LiteralText "E29" "  -1.58%  "
Set-LiteralText "D30" "555.80"
Set-LiteralText "E30" "  -3.63%  "
Set-LiteralText "D31" "8.04"
Set-LiteralText "E31" "  +0.69%  "
Set-LiteralText "E32" "  -0.75%  "
Set-LiteralText "E33" "  -0.60%  "
Set-LiteralText "E35" "  -2.11%  "
Set-LiteralText "E36" "  -0.14%  "
Set-LiteralText "E37" "  +1.34%  "
Set-LiteralText "D38" "19.71"
Set-LiteralText "E38" "  +2.21%  "
Set-LiteralText "E39" "  +0.96%  "
Set-LiteralText "E40" "  -3.37%  "
Set-LiteralText "E41" "  -1.72%  "
Set-LiteralText "E44" "  -2.03%  "
Set-LiteralText "D46" "158.71"
Set-LiteralText "E46" "  +1.40%  "
Set-LiteralText "E47" "  -0.46%  "
Set-LiteralText "D48" "22.03"
Set-LiteralText "E48" "  +0.00%  "
Set-LiteralText "E49" "  +0.12%  "
Set-LiteralText "E50" "  -2.22%  "
Set-LiteralText "E51" "  -0.91%  "

# Clean up the scratch cell entirely (value + formatting) so it leaves no trace.
$ws.Range("Z1").Clear()
$excel.CutCopyMode = $false
